$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 20 de Octubre de 2020 a las 23:46"

# Row 4
$ws.Range("B4").Value = 8510330
$ws.Range("C4").Value = 53677
$ws.Range("D4").Value = 5537355
$ws.Range("E4").Value = 2746961
$ws.Range("G4").Value = 792
$ws.Range("H4").Value = 226014

# Row 5
$ws.Range("B5").Value = 7648920
$ws.Range("C5").Value = 54184
$ws.Range("D5").Value = 6792334
$ws.Range("E5").Value = 740638
$ws.Range("G5").Value = 712
$ws.Range("H5").Value = 115948

# Row 6
$ws.Range("B6").Value = 5273954
$ws.Range("C6").Value = 22827
$ws.Range("E6").Value = 437458
$ws.Range("G6").Value = 611
$ws.Range("H6").Value = 154837

# Row 58
$ws.Range("B58").Value = 78533
$ws.Range("C58").Value = 309
$ws.Range("D58").Value = 75089
$ws.Range("E58").Value = 3139
$ws.Range("G58").Value = 3
$ws.Range("H58").Value = 305

# Row 84
$ws.Range("A84").Value = "Bulgaria"
$ws.Range("B84").Value = 31863
$ws.Range("C84").Value = 1336
$ws.Range("D84").Value = 17414
$ws.Range("E84").Value = 13430
$ws.Range("G84").Value = 11
$ws.Range("H84").Value = 1019

# Row 85
$ws.Range("A85").Value = "Eslovaquia"
$ws.Range("B85").Value = 31400
$ws.Range("C85").Value = 705
$ws.Range("D85").Value = 8004
$ws.Range("E85").Value = 23298
$ws.Range("G85").Value = 6
$ws.Range("H85").Value = 98

# Row 93
$ws.Range("B93").Value = 20342
$ws.Range("C93").Value = 18
$ws.Range("D93").Value = 20044
$ws.Range("E93").Value = 177

# Row 106
$ws.Range("A106").Value = "Maldivas"
$ws.Range("B106").Value = 11271
$ws.Range("C106").Value = 39
$ws.Range("D106").Value = 10234
$ws.Range("E106").Value = 1000
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 37

# Row 107
$ws.Range("A107").Value = "Luxemburgo"
$ws.Range("B107").Value = 11241
$ws.Range("C107").Value = 231
$ws.Range("D107").Value = 8471
$ws.Range("E107").Value = 2634
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 136

# Row 116
$ws.Range("B116").Value = 8187
$ws.Range("C116").Value = 28
$ws.Range("D116").Value = 7692
$ws.Range("E116").Value = 262
$ws.Range("G116").Value = 1
$ws.Range("H116").Value = 233

# Row 146
$ws.Range("B146").Value = 3796
$ws.Range("C146").Value = 31
$ws.Range("D146").Value = 2796
$ws.Range("E146").Value = 886
$ws.Range("G146").Value = 3
$ws.Range("H146").Value = 114

# Row 161
$ws.Range("B161").Value = 2104
$ws.Range("C161").Value = 33
$ws.Range("D161").Value = 1555
$ws.Range("E161").Value = 498

# Row 193
$ws.Range("B193").Value = 188
$ws.Range("C193").Value = 3
$ws.Range("D193").Value = 174
$ws.Range("E193").Value = 5

# Row 215
$ws.Range("D215").Value = 16
$ws.Range("E215").Value = 0
